$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.871.30"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "1.769.74"
$ws.Range("E3").Value = "  +1.07%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'327.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").Value = "'0.4485"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.22%  "
$ws.Range("D8").Value = "'0.3556"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").Value = "'0.07460"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "'42.06"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.15%  "
$ws.Range("D11").Value = "'1.099"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "'20.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "'6.034"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").Value = "'7.237"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.35%  "
$ws.Range("D16").Value = "1.769.38"
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("D17").Value = "'93.31"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").Value = "'0.06431"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("E21").Value = "  +2.85%  "
$ws.Range("D22").Value = "'5.782"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").Value = "27.917.92"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").Value = "'11.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("D25").Value = "'2.109"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("D26").Value = "'162.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("D28").Value = "1.973.61"
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("D29").Value = "'2.163"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.91%  "
$ws.Range("D30").Value = "'125.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").Value = "'1.095"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.86%  "
$ws.Range("D32").Value = "'0.09167"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").Value = "'5.594"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("D34").Value = "'3.642"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.75%  "
$ws.Range("D35").Value = "'11.88"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.44%  "
$ws.Range("D36").Value = "'0.02298"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("D37").Value = "'0.06122"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.6314"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "'4.969"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").Value = "'1.182"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.76%  "
$ws.Range("E42").Value = "  +0.92%  "
$ws.Range("D43").Value = "'7.923"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.28%  "
$ws.Range("D44").Value = "'13.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").Value = "'3.742"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("D46").Value = "'0.5871"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("D47").Value = "'122.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").Value = "'1.952"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("D49").Value = "'0.06903"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("D50").Value = "'1.137"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("E51").Value = "  +1.66%  "
